# Gene family file parser updates - species name converter sheet fixes.
# Strip the stray apostrophes from the two common-name entries that were
# causing parsing issues ("ma's_night_monkey" -> "mas_night_monkey" and
# "Hoffmann's_two_fingered_sloth" -> "Hoffmanns_two_fingered_sloth"),
# then leave the selection where the edits were made.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 26 = aotus_nancymaae / ma's_night_monkey / ENSANA
$ws.Range("B26").Value2 = "mas_night_monkey"

# Row 49 = choloepus_hoffmanni / Hoffmann's_two_fingered_sloth / ENSCHO
$ws.Range("B49").Value2 = "Hoffmanns_two_fingered_sloth"

# Reflect the scrolled/selected state from the session that made the edits.
$win = $excel.ActiveWindow
$win.ScrollRow = 27
$win.ScrollColumn = 1
$ws.Range("N46").Select()
